# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column headers in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, border, centered) from an existing
# header cell onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every data row (rows 2-47): every team/row
# shares the same season record (54 wins, 61 losses, 0 ties).
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 54
    $ws.Cells.Item($r, 31).Value = 61
    $ws.Cells.Item($r, 32).Value = 0
}
